# Fill in Andrew's section of the costing summary workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Andrew")

# Week 1
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2.5

# Week 2
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 2

# Week 3
$ws.Range("D4").Value = 9.5

# Week 4
$ws.Range("B5").Value = 0.5
$ws.Range("D5").Value = 6

# Week 5
$ws.Range("B6").Value = 4

# Week 6
$ws.Range("B7").Value = 0.5
$ws.Range("C7").Value = 2

# Week 7
$ws.Range("C8").Value = 1.5
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 1

# Week 8
$ws.Range("C9").Value = 2
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 1

# Week 9
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 2

# Week 10
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 1

# Week 11
$ws.Range("E12").Value = 3.5
$ws.Range("F12").Value = 0.5

# Greg's sheet was previously the selected tab with D26 selected; move its
# selection back to B1 and leave it deselected as a tab.
$wsGreg = $wb.Worksheets.Item("Greg")
$wsGreg.Activate()
$wsGreg.Range("B1").Select()

# Andrew becomes the active/selected tab, with G12 as the active cell.
$ws.Activate()
$ws.Range("G12").Select()
